$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 78900.96623825973
$ws.Range("C2").Value = 126625.1774154419
$ws.Range("D2").Value = 179616.2974086792
$ws.Range("E2").Value = 219277.3863073432
$ws.Range("B3").Value = 86790.95600780808
$ws.Range("C3").Value = 141834.2692948607
$ws.Range("D3").Value = 203230.6592846719
$ws.Range("E3").Value = 253024.8626478237
$ws.Range("B4").Value = 47473.32145807893
$ws.Range("C4").Value = 82650.45962447385
$ws.Range("D4").Value = 136052.5618903024
$ws.Range("E4").Value = 183271.8678927773
$ws.Range("B5").Value = 20985.77995440275
$ws.Range("C5").Value = 35479.66906454174
$ws.Range("D5").Value = 50568.11281908691
$ws.Range("E5").Value = 60392.30256534647
$ws.Range("B6").Value = 9508.61647798846
$ws.Range("C6").Value = 14196.46344123902
$ws.Range("D6").Value = 22621.75638317832
$ws.Range("E6").Value = 26644.84328347797
$ws.Range("B7").Value = 1671.270963587677
$ws.Range("C7").Value = 2809.121889588123
$ws.Range("D7").Value = 4308.994386236351
$ws.Range("E7").Value = 5118.22199367451
$ws.Range("B8").Value = 89345.77438605897
$ws.Range("C8").Value = 164238.2231598518
$ws.Range("D8").Value = 227769.7944633954
$ws.Range("E8").Value = 268666.5082649724
$ws.Range("B9").Value = 36187.57339999937
$ws.Range("C9").Value = 59721.89583435358
$ws.Range("D9").Value = 79895.1403610858
$ws.Range("E9").Value = 88557.3731103122
$ws.Range("B10").Value = 23147.54716588733
$ws.Range("C10").Value = 40216.96354845799
$ws.Range("D10").Value = 65042.17812641605
$ws.Range("E10").Value = 86793.06202533137
$ws.Range("B11").Value = 3719.188249141922
$ws.Range("C11").Value = 6050.047122912442
$ws.Range("D11").Value = 9957.26167774932
$ws.Range("E11").Value = 14669.66185150074
$ws.Range("B12").Value = 2194.418494361083
$ws.Range("C12").Value = 4528.928167216951
$ws.Range("D12").Value = 11320.9165339639
$ws.Range("E12").Value = 18606.31411736884
$ws.Range("B13").Value = 9949.239781516966
$ws.Range("C13").Value = 17391.14458628593
$ws.Range("D13").Value = 30564.50792310344
$ws.Range("E13").Value = 43170.97631742917
